$d = $word.ActiveDocument
$bullet = [char]0x2022

# --- Change 1: Collapse the three detailed "CORE COMPETENCIES" paragraphs
#     into a single summary paragraph with just the three headers.
$d.Paragraphs.Item(6).Range.Text = "Survey Methodology & Research Design $bullet Redistricting & Geospatial Analysis $bullet Data Analysis & Visualization"
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

# --- Change 2: Insert a new "TECHNICAL SKILLS" section (heading + three
#     detail paragraphs) right before the closing "For a more detailed..."
#     paragraph, which is now the last paragraph in the document.
$insertIndex = $d.Paragraphs.Count

$d.Paragraphs.Item($insertIndex).Range.InsertParagraphBefore()
$d.Paragraphs.Item($insertIndex).Range.Text = "TECHNICAL SKILLS"
$d.Paragraphs.Item($insertIndex).Style = "Heading2"
$insertIndex = $insertIndex + 1

$d.Paragraphs.Item($insertIndex).Range.InsertParagraphBefore()
$d.Paragraphs.Item($insertIndex).Range.Text = "SURVEY METHODOLOGY & RESEARCH DESIGN Survey Design and Questionnaire Development for Political and Market Research; Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR); Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling; Expert Testimony and Consultation on Research Methodology"
$insertIndex = $insertIndex + 1

$d.Paragraphs.Item($insertIndex).Range.InsertParagraphBefore()
$d.Paragraphs.Item($insertIndex).Range.Text = "REDISTRICTING & GEOSPATIAL ANALYSIS Redistricting Software Development and Boundary Estimation Systems; Geospatial Analysis; Choropleths and Hexagonal Grid Maps for Demographic Visualization; Court Case Analysis and Expert Testimony for Redistricting"
$insertIndex = $insertIndex + 1

$d.Paragraphs.Item($insertIndex).Range.InsertParagraphBefore()
$d.Paragraphs.Item($insertIndex).Range.Text = "DATA ANALYSIS & VISUALIZATION Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation); Data Visualization; Consumer Behavior Analysis and Market Segmentation; Multi-million Dollar Research Project Management"
$insertIndex = $insertIndex + 1
